# Fill in the missing clock-in/clock-out times for 6/11 and 6/12 (rows 39-40)
# on the Timesheet. G39/G40/G45 are formulas and will recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D39").Value = 0.50347222222222221
$ws.Range("E39").Value = 0.52361111111111114
$ws.Range("F39").Value = 0.61319444444444449

$ws.Range("C40").Value = 0.38611111111111113
$ws.Range("D40").Value = 0.52430555555555558
$ws.Range("E40").Value = 0.53680555555555554
$ws.Range("F40").Value = 0.61458333333333337

# Move the active selection to F36 on Sheet1 (matches the saved view state).
$ws.Activate()
$ws.Range("F36").Select()
